$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue 'D2' '28.357.00'
Set-TextValue 'E2' '  -0.59%  '

Set-TextValue 'D3' '1.568.75'
Set-TextValue 'E3' '  +0.25%  '

Set-TextValue 'E4' '  +0.00%  '

Set-TextValue 'D5' '210.94'
Set-TextValue 'E5' '  -0.35%  '

Set-TextValue 'E6' '  -0.55%  '

Set-TextValue 'E7' '  +0.02%  '

Set-TextValue 'D8' '44.26'
Set-TextValue 'E8' '  -4.07%  '

Set-TextValue 'D9' '23.82'
Set-TextValue 'E9' '  -1.50%  '

Set-TextValue 'E10' '  -0.94%  '

Set-TextValue 'E11' '  -0.77%  '

Set-TextValue 'D12' '0.0895'
Set-TextValue 'E12' '  +1.39%  '

Set-TextValue 'E13' '  +0.33%  '

Set-TextValue 'D14' '1.567.18'
Set-TextValue 'E14' '  +0.20%  '

Set-TextValue 'D15' '3.67'
Set-TextValue 'E15' '  -0.33%  '

Set-TextValue 'D16' '28.336.76'
Set-TextValue 'E16' '  -0.65%  '

Set-TextValue 'E17' '  -1.27%  '

Set-TextValue 'D18' '61.02'
Set-TextValue 'E18' '  -1.48%  '

Set-TextValue 'D19' '227.43'
Set-TextValue 'E19' '  +0.05%  '

Set-TextValue 'D20' '7.38'
Set-TextValue 'E20' '  +0.65%  '

Set-TextValue 'D21' '0.0₃0681'
Set-TextValue 'E21' '  -1.85%  '

Set-TextValue 'E22' '  +0.00%  '

Set-TextValue 'D23' '3.94'
Set-TextValue 'E23' '  +1.75%  '

Set-TextValue 'D24' '8.96'
Set-TextValue 'E24' '  -1.99%  '

Set-TextValue 'E25' '  -0.45%  '

Set-TextValue 'D26' '150.77'
Set-TextValue 'E26' '  +0.21%  '

Set-TextValue 'D27' '14.89'
Set-TextValue 'E27' '  -0.43%  '

Set-TextValue 'E28' '  -0.52%  '

Set-TextValue 'E29' '  -1.58%  '

Set-TextValue 'E30' '  +0.02%  '

Set-TextValue 'E31' '  +3.19%  '

Set-TextValue 'E32' '  -2.71%  '

Set-TextValue 'E33' '  -0.83%  '

Set-TextValue 'D34' '3.08'
Set-TextValue 'E34' '  -2.01%  '

Set-TextValue 'D35' '1.378.09'
Set-TextValue 'E35' '  -1.23%  '

Set-TextValue 'E36' '  +2.01%  '

Set-TextValue 'E37' '  -3.05%  '

Set-TextValue 'E38' '  -0.29%  '

Set-TextValue 'E39' '  +2.77%  '

Set-TextValue 'E40' '  -2.04%  '

Set-TextValue 'D41' '0.520'
Set-TextValue 'E41' '  -2.96%  '

Set-TextValue 'E42' '  +3.78%  '

Set-TextValue 'E43' '  +0.03%  '

Set-TextValue 'E44' '  -0.42%  '

Set-TextValue 'D45' '0.781'
Set-TextValue 'E45' '  -0.82%  '

Set-TextValue 'E46' '  -3.69%  '

Set-TextValue 'D47' '62.27'
Set-TextValue 'E47' '  -0.89%  '

Set-TextValue 'B48' 'RocketPoolETH'
Set-TextValue 'C48' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D48' '1.706.02'
Set-TextValue 'E48' '  +0.36%  '

Set-TextValue 'B49' 'WEMIXToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D49' '0.916'
Set-TextValue 'E49' '  -6.40%  '

Set-TextValue 'D50' '85.35'
Set-TextValue 'E50' '  -0.94%  '

Set-TextValue 'E51' '  -0.81%  '
